# Rename the existing sheet and add a new "Coordenadas" sheet with projected
# coordinate data, following the layout used for "Proyecciones" (old Sheet1).

$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> Proyecciones ---------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Proyecciones"

# Update header labels "Proy N"/"Proy E" -> "Proy Y"/"Proy X"
$ws1.Range("H1").Value = "Proy Y"
$ws1.Range("I1").Value = "Proy X"

# --- Add the new "Coordenadas" sheet right after "Proyecciones" ----------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Coordenadas"

# Header row
$ws2.Range("A1").Value = "Proy Y"
$ws2.Range("B1").Value = "Proy X"
$ws2.Range("C1").Value = "Corr Y"
$ws2.Range("D1").Value = "Corr X"
$ws2.Range("E1").Value = "Coord N"
$ws2.Range("F1").Value = "Coord E"

$headerRange = $ws2.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Blank separator row
$ws2.Range("A2:F2").Value = " "

# Data block 1
$ws2.Range("A3").Value = -40.192
$ws2.Range("B3").Value = -58.299
$ws2.Range("C3").Value = 0.004
$ws2.Range("D3").Value = -0.004
$ws2.Range("E3").Value = 1115.933
$ws2.Range("F3").Value = 2161.421

$ws2.Range("A4:F4").Value = " "

# Data block 2
$ws2.Range("A5").Value = 60.445
$ws2.Range("B5").Value = -25.621
$ws2.Range("C5").Value = 0.004
$ws2.Range("D5").Value = -0.004
$ws2.Range("E5").Value = 1075.746
$ws2.Range("F5").Value = 2103.118

$ws2.Range("A6:F6").Value = " "

# Data block 3
$ws2.Range("A7").Value = -20.267
$ws2.Range("B7").Value = 83.933
$ws2.Range("C7").Value = 0.005
$ws2.Range("D7").Value = -0.005
$ws2.Range("E7").Value = 1136.195
$ws2.Range("F7").Value = 2077.493

$ws2.Range("A8:F8").Value = " "

# Final row closes the loop back to the first coordinate
$ws2.Range("A9:D9").Value = " "
$ws2.Range("E9").Value = 1115.933
$ws2.Range("F9").Value = 2161.421
